# ---------------------------------------------------------------------------
# correct-questionnaire-scenario.pptx -- apply the authored edit:
#   1. Refresh the cached "datetimeFigureOut" field text (slide master and
#      every slide layout) from 27/11/15 to 03/02/16.
#   2. Re-word the caption of part (e) in the big screenshot collage: split
#      the trailing sentence into several runs and fix the wording from
#      "question 11 &12" to "questions 11 & 12".
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder refresh.
# ---------------------------------------------------------------------------
$newDate = "03/02/16"
$m = $p.SlideMaster

for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.Name -like "*Date*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $m.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.Name -like "*Date*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Caption (e) re-wording.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item(1)
$capTextBox = $group.GroupItems.Item("TextBox 14")
$tr = $capTextBox.TextFrame.TextRange

# The shape auto-sizes to its text ("resize shape to fit text"); remember
# the original height (in EMU, converted to points) so it can be restored
# once the wording change is done -- the caption box itself doesn't move.
$origHeightPt = 1077219 / 12700

# Replace the whole descriptive run (everything after the bold "(e)") with
# the corrected wording, as a single run first.
$tr.Characters(4, 108).Text = " deselecting question 8 deactivates the first group, such that the original questions 11 & 12 are active again"

# Now break it up into the individual runs of the authored edit, working
# from the end of the range back towards the start so earlier offsets are
# not invalidated by the replacements that follow.
$tr.Characters(98, 16).Text = "are active again"
$tr.Characters(93, 5).Text = "& 12 "
$tr.Characters(90, 3).Text = "11 "
$tr.Characters(80, 10).Text = "questions "

# Restore the caption box's original size (auto-fit would otherwise shrink
# it to match the new text).
$capTextBox.Height = $origHeightPt
